# Generate Report for Handback
#
# The localization-status workbook tracks two files that are going through
# hand-back: fc8743c1-6b31-4995-87bc-73c1dd8b615a.md and
# 4db536cf-0e30-4199-86c7-35de2277cff6.md.
#
# A fresh report run recorded 4db536cf's hand-back (it moved from "Ready for
# handoff" to "Handed back: in sync with en-US" with new handoff/handback
# timestamps), and at the same time swapped which row each file's data landed
# on (4db536cf now occupies row 2, fc8743c1 now occupies row 3) on every
# sheet. The per-row hyperlink targets (r:id) stay anchored to their original
# rows/URLs — only the visible link text and cell values move.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 (was fc8743c1, now 4db536cf) / Row 3 (was 4db536cf, now fc8743c1)
$ov.Range("A2").Value = "4db536cf-0e30-4199-86c7-35de2277cff6.md"
$ov.Range("B2").Value = "e2e\4db536cf-0e30-4199-86c7-35de2277cff6.md"
$ov.Range("G2").Value = "2016-10-21 03:56:28"

$ov.Range("A3").Value = "fc8743c1-6b31-4995-87bc-73c1dd8b615a.md"
$ov.Range("B3").Value = "e2e\fc8743c1-6b31-4995-87bc-73c1dd8b615a.md"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Range("G3").Value = "2016-10-21 03:54:28"

# Hyperlinks keep their original targets (r:id -> URL mapping unchanged);
# only the displayed text swaps between the two rows.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3c1074cac197a278493369710cb6d34e66f4d26a/e2e/fc8743c1-6b31-4995-87bc-73c1dd8b615a.md", "", "", "e2e\4db536cf-0e30-4199-86c7-35de2277cff6.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afbe306287bd6cff9a8c1d2dcbbcf8b9b8446c5a/e2e/4db536cf-0e30-4199-86c7-35de2277cff6.md", "", "", "e2e\fc8743c1-6b31-4995-87bc-73c1dd8b615a.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "4db536cf-0e30-4199-86c7-35de2277cff6.md"
$zh.Range("G2").Value = "4db536cf-0e30-4199-86c7-35de2277cff6.a6dd0ad62257a76d7082f828f8a65b56fdd00c62.zh-cn.xlf"
$zh.Range("H2").Value = "2016-10-21 03:56:17"
$zh.Range("I2").Value = "4db536cf-0e30-4199-86c7-35de2277cff6.md"
$zh.Range("J2").Value = "4db536cf-0e30-4199-86c7-35de2277cff6.a6dd0ad62257a76d7082f828f8a65b56fdd00c62.zh-cn.xlf"
$zh.Range("K2").Value = "2016-10-21 03:56:57"

$zh.Range("A3").Value = "fc8743c1-6b31-4995-87bc-73c1dd8b615a.md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("G3").Value = "fc8743c1-6b31-4995-87bc-73c1dd8b615a.7a0050f82e34a1476a14665c932ced0af1e49f4f.zh-cn.xlf"
$zh.Range("H3").Value = "2016-10-21 03:53:58"
$zh.Range("I3").Value = "fc8743c1-6b31-4995-87bc-73c1dd8b615a.md"
$zh.Range("J3").Value = "fc8743c1-6b31-4995-87bc-73c1dd8b615a.7a0050f82e34a1476a14665c932ced0af1e49f4f.zh-cn.xlf"
$zh.Range("P3").ClearContents()

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3c1074cac197a278493369710cb6d34e66f4d26a/e2e/fc8743c1-6b31-4995-87bc-73c1dd8b615a.md", "", "", "4db536cf-0e30-4199-86c7-35de2277cff6.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3853aa2d2b38fdc7e3d8d47557b5dd7eae33fc27/e2e/fc8743c1-6b31-4995-87bc-73c1dd8b615a.md", "", "", "4db536cf-0e30-4199-86c7-35de2277cff6.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afbe306287bd6cff9a8c1d2dcbbcf8b9b8446c5a/e2e/4db536cf-0e30-4199-86c7-35de2277cff6.md", "", "", "fc8743c1-6b31-4995-87bc-73c1dd8b615a.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3853aa2d2b38fdc7e3d8d47557b5dd7eae33fc27/e2e/4db536cf-0e30-4199-86c7-35de2277cff6.md", "", "", "fc8743c1-6b31-4995-87bc-73c1dd8b615a.md") | Out-Null

$zh.Columns.Item(16).ColumnWidth = 12.913719813028965

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "4db536cf-0e30-4199-86c7-35de2277cff6.md"
$de.Range("G2").Value = "4db536cf-0e30-4199-86c7-35de2277cff6.a6dd0ad62257a76d7082f828f8a65b56fdd00c62.de-de.xlf"
$de.Range("H2").Value = "2016-10-21 03:56:28"
$de.Range("I2").Value = "4db536cf-0e30-4199-86c7-35de2277cff6.md"
$de.Range("J2").Value = "4db536cf-0e30-4199-86c7-35de2277cff6.a6dd0ad62257a76d7082f828f8a65b56fdd00c62.de-de.xlf"
$de.Range("K2").Value = "2016-10-21 03:57:15"

$de.Range("A3").Value = "fc8743c1-6b31-4995-87bc-73c1dd8b615a.md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("G3").Value = "fc8743c1-6b31-4995-87bc-73c1dd8b615a.7a0050f82e34a1476a14665c932ced0af1e49f4f.de-de.xlf"
$de.Range("H3").Value = "2016-10-21 03:54:28"
$de.Range("I3").Value = "fc8743c1-6b31-4995-87bc-73c1dd8b615a.md"
$de.Range("J3").Value = "fc8743c1-6b31-4995-87bc-73c1dd8b615a.7a0050f82e34a1476a14665c932ced0af1e49f4f.de-de.xlf"
$de.Range("P3").ClearContents()

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3c1074cac197a278493369710cb6d34e66f4d26a/e2e/fc8743c1-6b31-4995-87bc-73c1dd8b615a.md", "", "", "4db536cf-0e30-4199-86c7-35de2277cff6.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/821e5bcc5adad1dbf71a4ca8adcb1b8682260dc4/e2e/fc8743c1-6b31-4995-87bc-73c1dd8b615a.md", "", "", "4db536cf-0e30-4199-86c7-35de2277cff6.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/afbe306287bd6cff9a8c1d2dcbbcf8b9b8446c5a/e2e/4db536cf-0e30-4199-86c7-35de2277cff6.md", "", "", "fc8743c1-6b31-4995-87bc-73c1dd8b615a.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/821e5bcc5adad1dbf71a4ca8adcb1b8682260dc4/e2e/4db536cf-0e30-4199-86c7-35de2277cff6.md", "", "", "fc8743c1-6b31-4995-87bc-73c1dd8b615a.md") | Out-Null

$de.Columns.Item(16).ColumnWidth = 12.913719813028965

Write-Output "Report regenerated for handback."
